# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.689.94"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "'3.178.82"
$ws.Range("E3").Value = "  -4.57%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'571.71"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "'171.66"
$ws.Range("E6").Value = "  -3.08%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.598"
$ws.Range("E8").Value = "  -2.77%  "

$ws.Range("D9").Value = "'3.177.11"
$ws.Range("E9").Value = "  -4.52%  "

$ws.Range("E10").Value = "  -1.86%  "

$ws.Range("D11").Value = "'6.62"
$ws.Range("E11").Value = "  -3.17%  "

$ws.Range("D12").Value = "'0.394"
$ws.Range("E12").Value = "  -2.94%  "

$ws.Range("D13").Value = "'3.727.28"
$ws.Range("E13").Value = "  -4.61%  "

$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("D15").Value = "'27.25"
$ws.Range("E15").Value = "  -3.73%  "

$ws.Range("D16").Value = "'65.642.08"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("E17").Value = "  -1.83%  "

$ws.Range("D18").Value = "'3.176.92"
$ws.Range("E18").Value = "  -4.81%  "

$ws.Range("D19").Value = "'5.75"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").Value = "'12.93"

$ws.Range("D21").Value = "'362.04"
$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("D22").Value = "'7.28"
$ws.Range("E22").Value = "  -1.66%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'68.80"
$ws.Range("E24").Value = "  -3.07%  "

$ws.Range("D25").Value = "'0.497"
$ws.Range("E25").Value = "  -3.62%  "

$ws.Range("D26").Value = "'3.308.34"

$ws.Range("E27").Value = "  -5.43%  "

$ws.Range("D28").Value = "'9.88"
$ws.Range("E28").Value = "  +3.95%  "

$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("E32").Value = "  -0.15%  "

$ws.Range("D33").Value = "'5.40"
$ws.Range("E33").Value = "  -3.16%  "

$ws.Range("D34").Value = "'22.09"
$ws.Range("E34").Value = "  -3.30%  "

$ws.Range("D35").Value = "'6.63"
$ws.Range("E35").Value = "  -2.29%  "

$ws.Range("E36").Value = "  -0.68%  "

$ws.Range("D37").Value = "'161.94"
$ws.Range("E37").Value = "  +1.26%  "

$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("D39").Value = "'0.835"
$ws.Range("E39").Value = "  -1.10%  "

$ws.Range("E40").Value = "  +3.67%  "

$ws.Range("D41").Value = "'26.51"
$ws.Range("E41").Value = "  -3.36%  "

$ws.Range("D42").Value = "'2.51"
$ws.Range("E42").Value = "  +1.33%  "

$ws.Range("D43").Value = "'2.652.60"
$ws.Range("E43").Value = "  -1.82%  "

$ws.Range("D44").Value = "'6.13"
$ws.Range("E44").Value = "  -0.57%  "

$ws.Range("D45").Value = "'4.21"
$ws.Range("E45").Value = "  -1.16%  "

$ws.Range("D46").Value = "'39.82"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").Value = "'0.0661"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("D48").Value = "'327.23"
$ws.Range("E48").Value = "  -1.67%  "

$ws.Range("D49").Value = "'23.87"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("E50").Value = "  -1.10%  "

$ws.Range("E51").Value = "  -0.43%  "
